# Chosing the Findbugs Analyses.xlsx
# Commit: "added a column to mark the current state of the implementation of a checker"
#
# Semantic change: a new column is inserted at column B of sheet "ausgewählte
# Analysen" that marks the implementation status of a checker ("Implemented" /
# "Needs IMDF" / ...). This shifts the previous columns B..H one to the right
# (B->C, C->D, D->E, F->G, G->H, H->I) and two rows (20 and 21) get new
# content in the (now) B and D columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ausgewählte Analysen")

# Insert a new column before column B - this shifts everything right and
# keeps all existing formulas/values/styles intact, matching the diff.
$ws.Columns("B:B").Insert()

# Row 20 (checker id 73 - Se: Non-serializable class has a serializable
# inner class): the (till now missing) bug-pattern-id column gets a value,
# and the status is now "Implemented".
$ws.Range("D20").Value = "SE_BAD_FIELD_INNER_CLASS"

# Row 21 (checker id 177 - NP: Store of null value into field annotated
# NonNull): bug-pattern-id, plus status "Needs IMDF".
$ws.Range("D21").Value = "NP_STORE_INTO_NONNULL_FIELD"

$ws.Range("B20").Value = "Implemented"
$ws.Range("B21").Value = "Needs IMDF"
